$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2165354330708661
$ws.Range("C2").Value = 0.5275590551181102
$ws.Range("J2").Value = 0.01181102362204724
$ws.Range("P2").Value = 0.1692913385826772
$ws.Range("S2").Value = 0.07480314960629922
$ws.Range("C3").Value = 0.0364963503649635
$ws.Range("J3").Value = 0.0218978102189781
$ws.Range("P3").Value = 0.708029197080292
$ws.Range("S3").Value = 0.2335766423357664
$ws.Range("J4").Value = 0.08108108108108109
$ws.Range("P4").Value = 0.5675675675675675
$ws.Range("S4").Value = 0.3513513513513514
$ws.Range("B6").Value = 0.03482587064676617
$ws.Range("D6").Value = 0.009950248756218905
$ws.Range("F6").Value = 0.03482587064676617
$ws.Range("J6").Value = 0.2686567164179104
$ws.Range("O6").Value = 0.009950248756218905
$ws.Range("Q6").Value = 0.1492537313432836
$ws.Range("R6").Value = 0.07960199004975124
$ws.Range("S6").Value = 0.4129353233830846
$ws.Range("B7").Value = 0.1355932203389831
$ws.Range("D7").Value = 0.02259887005649718
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.07909604519774012
$ws.Range("J7").Value = 0.1355932203389831
$ws.Range("Q7").Value = 0.1525423728813559
$ws.Range("R7").Value = 0.0847457627118644
$ws.Range("S7").Value = 0.384180790960452
$ws.Range("B8").Value = 0.08390022675736962
$ws.Range("D8").Value = 0.018140589569161
$ws.Range("F8").Value = 0.05215419501133787
$ws.Range("J8").Value = 0.1405895691609977
$ws.Range("O8").Value = 0.01360544217687075
$ws.Range("Q8").Value = 0.1473922902494331
$ws.Range("R8").Value = 0.1156462585034014
$ws.Range("S8").Value = 0.4285714285714285
$ws.Range("B9").Value = 0.05853658536585366
$ws.Range("D9").Value = 0.02439024390243903
$ws.Range("F9").Value = 0.05365853658536585
$ws.Range("J9").Value = 0.1024390243902439
$ws.Range("O9").Value = 0.01463414634146342
$ws.Range("Q9").Value = 0.1365853658536585
$ws.Range("R9").Value = 0.1317073170731707
$ws.Range("S9").Value = 0.4780487804878049
$ws.Range("B10").Value = 0.096
$ws.Range("D10").Value = 0.0152
$ws.Range("F10").Value = 0.07199999999999999
$ws.Range("J10").Value = 0.1208
$ws.Range("O10").Value = 0.0128
$ws.Range("Q10").Value = 0.2064
$ws.Range("R10").Value = 0.1056
$ws.Range("S10").Value = 0.3712
$ws.Range("G11").Value = 0.1355311355311355
$ws.Range("J11").Value = 0.1062271062271062
$ws.Range("K11").Value = 0.1868131868131868
$ws.Range("L11").Value = 0.5641025641025641
$ws.Range("S11").Value = 0.007326007326007326
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1923076923076923
$ws.Range("K12").Value = 0.01282051282051282
$ws.Range("L12").Value = 0.01923076923076923
$ws.Range("S12").Value = 0.02564102564102564
$ws.Range("G13").Value = 0.7894736842105263
$ws.Range("J13").Value = 0.1842105263157895
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.02577319587628866
$ws.Range("H15").Value = 0.1649484536082474
$ws.Range("I15").Value = 0.05154639175257732
$ws.Range("J15").Value = 0.3969072164948453
$ws.Range("K15").Value = 0.06185567010309279
$ws.Range("M15").Value = 0.01030927835051546
$ws.Range("O15").Value = 0.04123711340206185
$ws.Range("S15").Value = 0.2474226804123711
$ws.Range("F16").Value = 0.03821656050955414
$ws.Range("H16").Value = 0.1528662420382166
$ws.Range("I16").Value = 0.07643312101910828
$ws.Range("J16").Value = 0.4076433121019108
$ws.Range("K16").Value = 0.09554140127388536
$ws.Range("M16").Value = 0.01910828025477707
$ws.Range("O16").Value = 0.05732484076433121
$ws.Range("S16").Value = 0.1528662420382166
$ws.Range("F17").Value = 0.009876543209876543
$ws.Range("H17").Value = 0.1975308641975309
$ws.Range("I17").Value = 0.08888888888888889
$ws.Range("J17").Value = 0.4049382716049383
$ws.Range("K17").Value = 0.0962962962962963
$ws.Range("M17").Value = 0.03209876543209877
$ws.Range("N17").Value = 0.002469135802469136
$ws.Range("O17").Value = 0.06419753086419754
$ws.Range("S17").Value = 0.1037037037037037
$ws.Range("F18").Value = 0.01234567901234568
$ws.Range("H18").Value = 0.242798353909465
$ws.Range("I18").Value = 0.09465020576131687
$ws.Range("J18").Value = 0.3662551440329218
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.01234567901234568
$ws.Range("O18").Value = 0.06995884773662552
$ws.Range("S18").Value = 0.09053497942386832
$ws.Range("F19").Value = 0.009685230024213076
$ws.Range("H19").Value = 0.2025827280064568
$ws.Range("I19").Value = 0.09927360774818401
$ws.Range("J19").Value = 0.384180790960452
$ws.Range("K19").Value = 0.1008878127522195
$ws.Range("M19").Value = 0.01452784503631961
$ws.Range("N19").Value = 0.0008071025020177562
$ws.Range("O19").Value = 0.06698950766747377
$ws.Range("S19").Value = 0.1210653753026634

Write-Host "Applied 107 cell updates"
